# Week 8 Meeting.pptx edit script
# 1. Reposition the hero picture on slide 1 (Picture 3 / id 4).
# 2. Refresh the cached "datetime1" field text on the slide master and
#    every slide layout's Date placeholder from 11/22/2021 -> 4/8/2022.

$p = $ppt.ActivePresentation

# --- 1. Move/resize the background picture on slide 1 -----------------
# msoPicture = 13
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $sh = $s1.Shapes.Item($i)
    if ($sh.Name -eq "Picture 3" -or $sh.Type -eq 13) {
        $sh.Left = 0
        $sh.Top = -6.835434070866142
    }
}

# --- 2. Update the cached date text on the master + every layout ------
# ppPlaceholderDate = 16 ; msoPlaceholder = 14
$newDate = "4/8/2022"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16 -and $sh.HasTextFrame) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16 -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}
